# Update "想去人数" (column F) counts across the four worksheets to reflect
# the latest scrape (gh-pages output generated at 456a3b4).
# Only column F numeric values change; everything else is left untouched.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 1313
$ws.Range("F3").Value  = 92
$ws.Range("F4").Value  = 183
$ws.Range("F5").Value  = 988
$ws.Range("F6").Value  = 31
$ws.Range("F7").Value  = 1052
$ws.Range("F9").Value  = 859
$ws.Range("F11").Value = 742
$ws.Range("F12").Value = 1436
$ws.Range("F13").Value = 1043
$ws.Range("F14").Value = 765
$ws.Range("F15").Value = 775
$ws.Range("F16").Value = 91
$ws.Range("F17").Value = 596
$ws.Range("F18").Value = 108
$ws.Range("F19").Value = 680
$ws.Range("F20").Value = 1288
$ws.Range("F21").Value = 195
$ws.Range("F24").Value = 5314
$ws.Range("F25").Value = 284
$ws.Range("F27").Value = 2463
$ws.Range("F28").Value = 5903
$ws.Range("F30").Value = 1018
$ws.Range("F31").Value = 602
$ws.Range("F32").Value = 67
$ws.Range("F36").Value = 55
$ws.Range("F38").Value = 705
$ws.Range("F44").Value = 6
$ws.Range("F45").Value = 17
$ws.Range("F47").Value = 98
$ws.Range("F49").Value = 25

# --- Sheet: 演出 (Performances) ---------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value  = 10
$ws.Range("F6").Value  = 75
$ws.Range("F8").Value  = 124
$ws.Range("F9").Value  = 487
$ws.Range("F14").Value = 680
$ws.Range("F15").Value = 19
$ws.Range("F16").Value = 788
$ws.Range("F21").Value = 19
$ws.Range("F25").Value = 1709
$ws.Range("F31").Value = 108
$ws.Range("F35").Value = 40
$ws.Range("F38").Value = 67
$ws.Range("F40").Value = 496
$ws.Range("F42").Value = 33
$ws.Range("F43").Value = 4

# --- Sheet: 本地生活 (Local life) -------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 671
$ws.Range("F5").Value = 775
$ws.Range("F6").Value = 426
$ws.Range("F7").Value = 244

# --- Sheet: 全部类型 (All categories, union of the above) -------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 671
$ws.Range("F4").Value  = 1313
$ws.Range("F5").Value  = 92
$ws.Range("F6").Value  = 426
$ws.Range("F7").Value  = 244
$ws.Range("F8").Value  = 244
$ws.Range("F9").Value  = 75
$ws.Range("F10").Value = 988
$ws.Range("F11").Value = 124
$ws.Range("F12").Value = 31
$ws.Range("F13").Value = 1052
$ws.Range("F15").Value = 859
$ws.Range("F17").Value = 487
$ws.Range("F18").Value = 742
$ws.Range("F19").Value = 1436
$ws.Range("F21").Value = 1043
$ws.Range("F22").Value = 765
$ws.Range("F24").Value = 775
$ws.Range("F25").Value = 91
$ws.Range("F26").Value = 596
$ws.Range("F27").Value = 680
$ws.Range("F28").Value = 1288
$ws.Range("F29").Value = 195
$ws.Range("F31").Value = 5314
$ws.Range("F32").Value = 284
$ws.Range("F34").Value = 2463
$ws.Range("F35").Value = 5903
$ws.Range("F36").Value = 1018
$ws.Range("F37").Value = 1709
$ws.Range("F38").Value = 602
$ws.Range("F39").Value = 67
$ws.Range("F40").Value = 1054
$ws.Range("F41").Value = 55
$ws.Range("F42").Value = 705
$ws.Range("F47").Value = 496
$ws.Range("F48").Value = 17
$ws.Range("F49").Value = 33
